$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.399.48'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.500.36'
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.52'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.50'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.62'
$ws.Range("E9").Value = '  +5.74%  '

$ws.Range("E10").Value = '  +1.05%  '

$ws.Range("E11").Value = '  +4.30%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.096.64'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("E13").Value = '  +0.69%  '

$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.501.24'
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.75'
$ws.Range("E16").Value = '  +2.29%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.396.00'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.81'
$ws.Range("E19").Value = '  +1.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.59'
$ws.Range("E20").Value = '  -0.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.41'
$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.583'
$ws.Range("E22").Value = '  +2.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.639.77'
$ws.Range("E23").Value = '  +0.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.43'
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  -0.81%  '

$ws.Range("E27").Value = '  +2.57%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.45'
$ws.Range("E28").Value = '  +0.37%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E30").Value = '  +1.91%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.49'
$ws.Range("E31").Value = '  -3.19%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.21'
$ws.Range("E32").Value = '  -0.61%  '

$ws.Range("E33").Value = '  +5.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.528.17'
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.44'
$ws.Range("E36").Value = '  +0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("E37").Value = '  +1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.98'
$ws.Range("E38").Value = '  +1.70%  '

$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.00'
$ws.Range("E40").Value = '  +2.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0793'
$ws.Range("E41").Value = '  +1.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.810'
$ws.Range("E42").Value = '  +0.70%  '

$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("E44").Value = '  +1.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.06'
$ws.Range("E45").Value = '  -1.63%  '

$ws.Range("E46").Value = '  -0.24%  '

$ws.Range("E47").Value = '  +0.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.928'
$ws.Range("E48").Value = '  +3.41%  '

$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.82'
$ws.Range("E49").Value = '  +1.21%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.411.57'
$ws.Range("E50").Value = '  -2.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0260'
$ws.Range("E51").Value = '  +0.26%  '

